# Update scripts with new TPM-derived NATMI results for the Rspo1-Lgr6
# ligand-receptor pair: the "Resolving-Mac" sending cluster is dropped and
# replaced by a new "ECs" target cluster (for both FAPs and
# Inflammatory-Mac senders), and all numeric columns (E:T) are refreshed
# with their recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = "FAPs"
$ws.Cells.Item(2, 2).Value2 = "Rspo1"
$ws.Cells.Item(2, 3).Value2 = "Lgr6"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 0.07032533333333334
$ws.Cells.Item(2, 8).Value2 = 0.210976
$ws.Cells.Item(2, 9).Value2 = 0.1886527232569993
$ws.Cells.Item(2, 10).Value2 = 0.1886527232569993
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 13).Value2 = 0.01182833333333333
$ws.Cells.Item(2, 14).Value2 = 0.035485
$ws.Cells.Item(2, 15).Value2 = 0.03045376408867423
$ws.Cells.Item(2, 16).Value2 = 0.03045376408867423
$ws.Cells.Item(2, 17).Value2 = 0.0008318314844444446
$ws.Cells.Item(2, 18).Value2 = 0.007486483360000001
$ws.Cells.Item(2, 19).Value2 = 0.005745185528754602
$ws.Cells.Item(2, 20).Value2 = 0.005745185528754602

$ws.Cells.Item(3, 1).Value2 = "FAPs"
$ws.Cells.Item(3, 2).Value2 = "Rspo1"
$ws.Cells.Item(3, 3).Value2 = "Lgr6"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 0.07032533333333334
$ws.Cells.Item(3, 8).Value2 = 0.210976
$ws.Cells.Item(3, 9).Value2 = 0.1886527232569993
$ws.Cells.Item(3, 10).Value2 = 0.1886527232569993
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 0.2280253333333333
$ws.Cells.Item(3, 14).Value2 = 0.684076
$ws.Cells.Item(3, 15).Value2 = 0.5870843771374921
$ws.Cells.Item(3, 16).Value2 = 0.5870843771374921
$ws.Cells.Item(3, 17).Value2 = 0.01603595757511111
$ws.Cells.Item(3, 18).Value2 = 0.144323618176
$ws.Cells.Item(3, 19).Value2 = 0.1107550665286271
$ws.Cells.Item(3, 20).Value2 = 0.1107550665286271

$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Rspo1"
$ws.Cells.Item(4, 3).Value2 = "Lgr6"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 0.07032533333333334
$ws.Cells.Item(4, 8).Value2 = 0.210976
$ws.Cells.Item(4, 9).Value2 = 0.1886527232569993
$ws.Cells.Item(4, 10).Value2 = 0.1886527232569993
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 0.1485493333333333
$ws.Cells.Item(4, 14).Value2 = 0.445648
$ws.Cells.Item(4, 15).Value2 = 0.3824618587738337
$ws.Cells.Item(4, 16).Value2 = 0.3824618587738337
$ws.Cells.Item(4, 17).Value2 = 0.01044678138311111
$ws.Cells.Item(4, 18).Value2 = 0.094021032448
$ws.Cells.Item(4, 19).Value2 = 0.0721524711996176
$ws.Cells.Item(4, 20).Value2 = 0.0721524711996176

$ws.Cells.Item(5, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(5, 2).Value2 = "Rspo1"
$ws.Cells.Item(5, 3).Value2 = "Lgr6"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 0.3024513333333333
$ws.Cells.Item(5, 8).Value2 = 0.907354
$ws.Cells.Item(5, 9).Value2 = 0.8113472767430007
$ws.Cells.Item(5, 10).Value2 = 0.8113472767430007
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(5, 13).Value2 = 0.01182833333333333
$ws.Cells.Item(5, 14).Value2 = 0.035485
$ws.Cells.Item(5, 15).Value2 = 0.03045376408867423
$ws.Cells.Item(5, 16).Value2 = 0.03045376408867423
$ws.Cells.Item(5, 17).Value2 = 0.003577495187777779
$ws.Cells.Item(5, 18).Value2 = 0.03219745669
$ws.Cells.Item(5, 19).Value2 = 0.02470857855991962
$ws.Cells.Item(5, 20).Value2 = 0.02470857855991962

$ws.Cells.Item(6, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(6, 2).Value2 = "Rspo1"
$ws.Cells.Item(6, 3).Value2 = "Lgr6"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 0.3024513333333333
$ws.Cells.Item(6, 8).Value2 = 0.907354
$ws.Cells.Item(6, 9).Value2 = 0.8113472767430007
$ws.Cells.Item(6, 10).Value2 = 0.8113472767430007
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 0.2280253333333333
$ws.Cells.Item(6, 14).Value2 = 0.684076
$ws.Cells.Item(6, 15).Value2 = 0.5870843771374921
$ws.Cells.Item(6, 16).Value2 = 0.5870843771374921
$ws.Cells.Item(6, 17).Value2 = 0.06896656610044445
$ws.Cells.Item(6, 18).Value2 = 0.620699094904
$ws.Cells.Item(6, 19).Value2 = 0.476329310608865
$ws.Cells.Item(6, 20).Value2 = 0.476329310608865

$ws.Cells.Item(7, 1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(7, 2).Value2 = "Rspo1"
$ws.Cells.Item(7, 3).Value2 = "Lgr6"
$ws.Cells.Item(7, 4).Value2 = "MuSCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 0.3024513333333333
$ws.Cells.Item(7, 8).Value2 = 0.907354
$ws.Cells.Item(7, 9).Value2 = 0.8113472767430007
$ws.Cells.Item(7, 10).Value2 = 0.8113472767430007
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 0.1485493333333333
$ws.Cells.Item(7, 14).Value2 = 0.445648
$ws.Cells.Item(7, 15).Value2 = 0.3824618587738337
$ws.Cells.Item(7, 16).Value2 = 0.3824618587738337
$ws.Cells.Item(7, 17).Value2 = 0.04492894393244445
$ws.Cells.Item(7, 18).Value2 = 0.404360495392
$ws.Cells.Item(7, 19).Value2 = 0.3103093875742161
$ws.Cells.Item(7, 20).Value2 = 0.3103093875742161

